# Revert "Edit a big error in the diagram"
#
# 1) Title slide text was mistakenly overwritten with "Doer List";
#    restore it to "AddressBook" + " \u2013 Level 4" (two runs).
# 2) A trailing blank slide (Blank layout) was deleted; re-add it as the
#    last slide in the deck.

$p = $ppt.ActivePresentation

# --- 1. Fix the title placeholder text on slide 1 ------------------------
$titleSlide = $p.Slides.Item(1)
$titleShape = $titleSlide.Shapes.Item(1)
for ($i = 1; $i -le $titleSlide.Shapes.Count; $i++) {
    $candidate = $titleSlide.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -eq "Doer List") {
        $titleShape = $candidate
    }
}
$tr = $titleShape.TextFrame.TextRange
$tr.Text = "AddressBook"
$null = $tr.InsertAfter(" " + [char]0x2013 + " Level 4")

# --- 2. Re-add the missing blank slide at the end of the deck ------------
$lastIndex = $p.Slides.Count
$null = $p.Slides.Item($lastIndex).Duplicate()

$newSlide = $p.Slides.Item($lastIndex + 1)
for ($i = $newSlide.Shapes.Count; $i -ge 1; $i--) {
    $newSlide.Shapes.Item($i).Delete()
}
